{"js": "// Replace the 25 table-cell practice problems with their updated values.\n// Each old value is unique in the document, so a direct search/replace\n// per pair is unambiguous.\nconst replacements = [\n  [\"57\u00f74=14, 1\", \"84\u00f75=16, 4\"],\n  [\"95\u00f73=31, 2\", \"94\u00f78=11, 6\"],\n  [\"69\u00f78=8, 5\", \"82\u00f75=16, 2\"],\n  [\"60\u00f77=8, 4\", \"16\u00f77=2, 2\"],\n  [\"70\u00f72=35, 0\", \"50\u00f73=16, 2\"],\n  [\"59\u00f73=19, 2\", \"30\u00f77=4, 2\"],\n  [\"93\u00f78=11, 5\", \"58\u00f76=9, 4\"],\n  [\"92\u00f79=10, 2\", \"33\u00f75=6, 3\"],\n  [\"81\u00f76=13, 3\", \"19\u00f73=6, 1\"],\n  [\"50\u00f74=12, 2\", \"45\u00f72=22, 1\"],\n  [\"14\u00f76=2, 2\", \"74\u00f79=8, 2\"],\n  [\"20\u00f79=2, 2\", \"64\u00f77=9, 1\"],\n  [\"70\u00f77=10, 0\", \"10\u00f79=1, 1\"],\n  [\"16\u00f78=2, 0\", \"75\u00f76=12, 3\"],\n  [\"66\u00f79=7, 3\", \"61\u00f74=15, 1\"],\n  [\"85\u00f78=10, 5\", \"80\u00f78=10, 0\"],\n  [\"10\u00f74=2, 2\", \"75\u00f79=8, 3\"],\n  [\"45\u00f78=5, 5\", \"14\u00f72=7, 0\"],\n  [\"99\u00f78=12, 3\", \"47\u00f73=15, 2\"],\n  [\"45\u00f73=15, 0\", \"25\u00f73=8, 1\"],\n  [\"62\u00f76=10, 2\", \"16\u00f79=1, 7\"],\n  [\"67\u00f73=22, 1\", \"41\u00f78=5, 1\"],\n  [\"45\u00f76=7, 3\", \"45\u00f78=5, 5\"],\n  [\"30\u00f79=3, 3\", \"15\u00f74=3, 3\"],\n  [\"55\u00f77=7, 6\", \"55\u00f79=6, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 table-cell practice problems with their updated values.\n# Each old value is unique in the document, so Find/Replace per pair is\n# unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"57\u00f74=14, 1\", \"84\u00f75=16, 4\"),\n    @(\"95\u00f73=31, 2\", \"94\u00f78=11, 6\"),\n    @(\"69\u00f78=8, 5\", \"82\u00f75=16, 2\"),\n    @(\"60\u00f77=8, 4\", \"16\u00f77=2, 2\"),\n    @(\"70\u00f72=35, 0\", \"50\u00f73=16, 2\"),\n    @(\"59\u00f73=19, 2\", \"30\u00f77=4, 2\"),\n    @(\"93\u00f78=11, 5\", \"58\u00f76=9, 4\"),\n    @(\"92\u00f79=10, 2\", \"33\u00f75=6, 3\"),\n    @(\"81\u00f76=13, 3\", \"19\u00f73=6, 1\"),\n    @(\"50\u00f74=12, 2\", \"45\u00f72=22, 1\"),\n    @(\"14\u00f76=2, 2\", \"74\u00f79=8, 2\"),\n    @(\"20\u00f79=2, 2\", \"64\u00f77=9, 1\"),\n    @(\"70\u00f77=10, 0\", \"10\u00f79=1, 1\"),\n    @(\"16\u00f78=2, 0\", \"75\u00f76=12, 3\"),\n    @(\"66\u00f79=7, 3\", \"61\u00f74=15, 1\"),\n    @(\"85\u00f78=10, 5\", \"80\u00f78=10, 0\"),\n    @(\"10\u00f74=2, 2\", \"75\u00f79=8, 3\"),\n    @(\"45\u00f78=5, 5\", \"14\u00f72=7, 0\"),\n    @(\"99\u00f78=12, 3\", \"47\u00f73=15, 2\"),\n    @(\"45\u00f73=15, 0\", \"25\u00f73=8, 1\"),\n    @(\"62\u00f76=10, 2\", \"16\u00f79=1, 7\"),\n    @(\"67\u00f73=22, 1\", \"41\u00f78=5, 1\"),\n    @(\"45\u00f76=7, 3\", \"45\u00f78=5, 5\"),\n    @(\"30\u00f79=3, 3\", \"15\u00f74=3, 3\"),\n    @(\"55\u00f77=7, 6\", \"55\u00f79=6, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute(\n        $oldText,    # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $newText,    # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n}\n"}
